# "culture collection を MIxS から再度削除 INSDC2017 での確認に基づく"
#
# The "culture_collection" field (column W, shared-string index 22) is being
# removed again. Deleting the column shifts every subsequent column (and its
# shared-string value) one slot to the left, but this runtime's
# Range/Columns.Delete does not re-anchor the cell comments (notes) that sit
# on row 15 - they stay pinned to their original cell address. So, before
# deleting the column, push each comment's text one cell to the left by hand
# (W<-X, X<-Y, ... AJ<-AK, AK<-AL), drop the now-duplicated trailing comment
# at AL15, and only then delete column W so the cell/string data lines up
# with the already-shifted comments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("W15").Comment.Text('Plasmids that have significance phenotypic consequence') | Out-Null
$ws.Range("X15").Comment.Text('Health or disease status of sample at time of collection') | Out-Null
$ws.Range("Y15").Comment.Text('The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".') | Out-Null
$ws.Range("Z15").Comment.Text('NCBI taxonomy ID of the host, e.g. 9606') | Out-Null
$ws.Range("AA15").Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.') | Out-Null
$ws.Range("AB15").Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html') | Out-Null
$ws.Range("AC15").Comment.Text('To what is the entity pathogenic') | Out-Null
$ws.Range("AD15").Comment.Text('Primary publication or genome report in the form of pubmed ID, DOI or URL') | Out-Null
$ws.Range("AE15").Comment.Text('Method or device employed for collecting sample') | Out-Null
$ws.Range("AF15").Comment.Text('Processing applied to the sample during or after isolation') | Out-Null
$ws.Range("AG15").Comment.Text('Amount or size of sample (volume, mass or area) that was collected') | Out-Null
$ws.Range("AH15").Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.') | Out-Null
$ws.Range("AI15").Comment.Text('Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier') | Out-Null
$ws.Range("AJ15").Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)') | Out-Null
$ws.Range("AK15").Comment.Text('Feeding position in food chain (eg., chemolithotroph)') | Out-Null

$ws.Range("AL15").Comment.Delete() | Out-Null

$ws.Columns.Item(23).Delete() | Out-Null
